# Major refactoring of code
# Changed code to separate logic from generation of output.
#
# This script corrects a handful of text strings (typography / wording
# fixes) on both worksheets, resets a stray custom font on two cells
# back to the default style, and updates the saved cell selection on
# each sheet.

$wb = $excel.ActiveWorkbook

$wsExternal = $wb.Worksheets.Item("External Service")
$wsNU       = $wb.Worksheets.Item("Northwestern University Service")

# ----------------------------------------------------------------------
# "External Service" sheet
# ----------------------------------------------------------------------

# Row 6 (Editor / Journal of Statistical Mechanics...): drop trailing "s"
$wsExternal.Range("E6").Value = "Journal of Statistical Mechanics: Theory and Experiment"

# Row 7 (Program Committee / 5th Int'l Workshop...): spell out "International"
$wsExternal.Range("E7").Value = "5th International Workshop on Biosignal Interpretation"
$wsExternal.Range("F7").Value = "International Federation for Medical and Biological Engineering, International Medical Informatics Association, IEEE Engineering in Medicine and Biology, Japan Society of Medical Electronics and Biomedical Engineering"

# Row 11 (Program Committee / Modeling and Simulation Workshop): replace em-dash with double hyphen
$wsExternal.Range("F11").Value = "Department of Homeland Security -- Science and Technology"

# Row 16 had a stray custom font (Times, size 11.5) applied to E16/F16; restore default style
$wsExternal.Range("E16").Style = "Normal"
$wsExternal.Range("F16").Style = "Normal"

# Update the remembered selection for this sheet
$wsExternal.Activate()
$wsExternal.Range("F8").Select()

# ----------------------------------------------------------------------
# "Northwestern University Service" sheet
# ----------------------------------------------------------------------

# Rows 2, 13, 14, 15 (Preceptor rows): "Comment" column filled in with "Graduate Program"
$wsNU.Range("E2").Value = "Graduate Program"
$wsNU.Range("E13").Value = "Graduate Program"
$wsNU.Range("E14").Value = "Graduate Program"
$wsNU.Range("E15").Value = "Graduate Program"

# Row 21: drop the apostrophe-s after "Northwestern"
$wsNU.Range("E21").Value = "One Northwestern Task Force on Integrated Enterprise"

# Rows 29, 30: drop curly quotes around "Big Data"
$wsNU.Range("E29").Value = "L. Dumas Domain Dinner on Big Data"
$wsNU.Range("E30").Value = "Faculty Workshop on Big Data"

# Update the remembered selection for this sheet and leave it the active tab
$wsNU.Activate()
$wsNU.Range("E22").Select()
